# Update "想去人数" (interest count) figures in the F column across the
# workbook's sheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value2  = 2752
$wsExhibit.Range("F4").Value2  = 1085
$wsExhibit.Range("F5").Value2  = 20003
$wsExhibit.Range("F8").Value2  = 758
$wsExhibit.Range("F14").Value2 = 74
$wsExhibit.Range("F15").Value2 = 383
$wsExhibit.Range("F19").Value2 = 216

# 演出 (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value2 = 296

# 本地生活 (Local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value2 = 656

# 全部类型 (All types - combined sheet)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value2  = 656
$wsAll.Range("F8").Value2  = 2752
$wsAll.Range("F9").Value2  = 1085
$wsAll.Range("F10").Value2 = 20003
$wsAll.Range("F15").Value2 = 296
$wsAll.Range("F17").Value2 = 758
$wsAll.Range("F25").Value2 = 74
$wsAll.Range("F28").Value2 = 383
$wsAll.Range("F36").Value2 = 216
